$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D2").Value = "2016-02-15 08:19:46"
$wsZhCn.Range("G2").Value = "2016-02-15 08:20:48"

$wsDeDe.Range("D2").Value = "2016-02-15 08:20:01"
$wsDeDe.Range("G2").Value = "2016-02-15 08:21:14"
